$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(88, 8).Value = 439504.38
$ws_ALC.Cells.Item(88, 9).Value = 1112139.6
$ws_ALC.Cells.Item(88, 10).Value = 11463.818
$ws_ALC.Cells.Item(88, 11).Value = 1112139.6
$ws_ALC.Cells.Item(88, 12).Value = 11463.818
$ws_ALC.Cells.Item(88, 13).Value = -1111733.6
$ws_ALC.Cells.Item(88, 14).Value = -12275.818

$ws_ALC.Cells.Item(91, 8).Value = 439504.38
$ws_ALC.Cells.Item(91, 9).Value = 1112139.6
$ws_ALC.Cells.Item(91, 10).Value = 11463.818
$ws_ALC.Cells.Item(91, 11).Value = 1112139.6
$ws_ALC.Cells.Item(91, 12).Value = 11463.818
$ws_ALC.Cells.Item(91, 13).Value = -1110735.6
$ws_ALC.Cells.Item(91, 14).Value = -14271.818

$ws_ALC.Cells.Item(123, 8).Value = 30778.889
$ws_ALC.Cells.Item(123, 10).Value = 30778.889
$ws_ALC.Cells.Item(123, 12).Value = 30778.889
$ws_ALC.Cells.Item(123, 14).Value = -40578.889

$ws_ALC.Cells.Item(137, 8).Value = 1456883
$ws_ALC.Cells.Item(137, 9).Value = 2854735.2
$ws_ALC.Cells.Item(137, 10).Value = 5267.154
$ws_ALC.Cells.Item(137, 11).Value = 8564205.600000001
$ws_ALC.Cells.Item(137, 12).Value = 15801.462
$ws_ALC.Cells.Item(137, 13).Value = -8561655.600000001
$ws_ALC.Cells.Item(137, 14).Value = -20901.462

$ws_ALC.Cells.Item(138, 8).Value = 2265.5334
$ws_ALC.Cells.Item(138, 9).Value = 1277.125
$ws_ALC.Cells.Item(138, 10).Value = 3395.1428
$ws_ALC.Cells.Item(138, 11).Value = 3831.375
$ws_ALC.Cells.Item(138, 12).Value = 10185.4284
$ws_ALC.Cells.Item(138, 13).Value = 1308.625
$ws_ALC.Cells.Item(138, 14).Value = -20465.4284

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(32, 8).Value = 10557.211
$ws_ARM.Cells.Item(32, 9).Value = 10083.667
$ws_ARM.Cells.Item(32, 11).Value = 10083.667
$ws_ARM.Cells.Item(32, 13).Value = -9796.666999999999

$ws_ARM.Cells.Item(61, 8).Value = 1929.2051
$ws_ARM.Cells.Item(61, 9).Value = 1288.1177
$ws_ARM.Cells.Item(61, 10).Value = 2424.5908
$ws_ARM.Cells.Item(61, 11).Value = 1288.1177
$ws_ARM.Cells.Item(61, 12).Value = 2424.5908
$ws_ARM.Cells.Item(61, 13).Value = -1076.1177
$ws_ARM.Cells.Item(61, 14).Value = -2848.5908

$ws_ARM.Cells.Item(136, 8).Value = 1929.2051
$ws_ARM.Cells.Item(136, 9).Value = 1288.1177
$ws_ARM.Cells.Item(136, 10).Value = 2424.5908
$ws_ARM.Cells.Item(136, 11).Value = 3864.3531
$ws_ARM.Cells.Item(136, 12).Value = 7273.7724
$ws_ARM.Cells.Item(136, 13).Value = -1314.3531
$ws_ARM.Cells.Item(136, 14).Value = -12373.7724

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(94, 8).Value = 1195.5333
$ws_BSM.Cells.Item(94, 9).Value = 1223.7858
$ws_BSM.Cells.Item(94, 10).Value = 800
$ws_BSM.Cells.Item(94, 11).Value = 1223.7858
$ws_BSM.Cells.Item(94, 12).Value = 800
$ws_BSM.Cells.Item(94, 13).Value = -772.7858000000001
$ws_BSM.Cells.Item(94, 14).Value = -1702

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(31, 8).Value = 1542.32
$ws_CRP.Cells.Item(31, 9).Value = 957.9583
$ws_CRP.Cells.Item(31, 10).Value = 2081.7307
$ws_CRP.Cells.Item(31, 11).Value = 957.9583
$ws_CRP.Cells.Item(31, 12).Value = 2081.7307
$ws_CRP.Cells.Item(31, 13).Value = -662.9583
$ws_CRP.Cells.Item(31, 14).Value = -2671.7307

$ws_CRP.Cells.Item(34, 8).Value = 1542.32
$ws_CRP.Cells.Item(34, 9).Value = 957.9583
$ws_CRP.Cells.Item(34, 10).Value = 2081.7307
$ws_CRP.Cells.Item(34, 11).Value = 957.9583
$ws_CRP.Cells.Item(34, 12).Value = 2081.7307
$ws_CRP.Cells.Item(34, 13).Value = -755.9583
$ws_CRP.Cells.Item(34, 14).Value = -2485.7307

$ws_CRP.Cells.Item(58, 8).Value = 2171.353
$ws_CRP.Cells.Item(58, 9).Value = 1274.9131
$ws_CRP.Cells.Item(58, 10).Value = 4045.7273
$ws_CRP.Cells.Item(58, 11).Value = 1274.9131
$ws_CRP.Cells.Item(58, 12).Value = 4045.7273
$ws_CRP.Cells.Item(58, 13).Value = -1071.9131
$ws_CRP.Cells.Item(58, 14).Value = -4451.7273

$ws_CRP.Cells.Item(136, 8).Value = 2171.353
$ws_CRP.Cells.Item(136, 9).Value = 1274.9131
$ws_CRP.Cells.Item(136, 10).Value = 4045.7273
$ws_CRP.Cells.Item(136, 11).Value = 3824.7393
$ws_CRP.Cells.Item(136, 12).Value = 12137.1819
$ws_CRP.Cells.Item(136, 13).Value = -1274.7393
$ws_CRP.Cells.Item(136, 14).Value = -17237.1819

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(68, 8).Value = 1344.3068
$ws_CUL.Cells.Item(68, 9).Value = 1043.5264
$ws_CUL.Cells.Item(68, 10).Value = 1427.1305
$ws_CUL.Cells.Item(68, 11).Value = 3130.5792
$ws_CUL.Cells.Item(68, 12).Value = 4281.3915
$ws_CUL.Cells.Item(68, 13).Value = -2319.5792
$ws_CUL.Cells.Item(68, 14).Value = -5903.3915

$ws_CUL.Cells.Item(71, 8).Value = 1344.3068
$ws_CUL.Cells.Item(71, 9).Value = 1043.5264
$ws_CUL.Cells.Item(71, 10).Value = 1427.1305
$ws_CUL.Cells.Item(71, 11).Value = 9391.7376
$ws_CUL.Cells.Item(71, 12).Value = 12844.1745
$ws_CUL.Cells.Item(71, 13).Value = -5335.7376
$ws_CUL.Cells.Item(71, 14).Value = -20956.1745

$ws_CUL.Cells.Item(107, 8).Value = 10671.904
$ws_CUL.Cells.Item(107, 9).Value = 10442.9
$ws_CUL.Cells.Item(107, 10).Value = 10880.091
$ws_CUL.Cells.Item(107, 11).Value = 31328.7
$ws_CUL.Cells.Item(107, 12).Value = 32640.273
$ws_CUL.Cells.Item(107, 13).Value = -29408.7
$ws_CUL.Cells.Item(107, 14).Value = -36480.273

$ws_CUL.Cells.Item(129, 8).Value = 98166.87
$ws_CUL.Cells.Item(129, 9).Value = 334104
$ws_CUL.Cells.Item(129, 10).Value = 1647.1364
$ws_CUL.Cells.Item(129, 11).Value = 1002312
$ws_CUL.Cells.Item(129, 12).Value = 4941.4092
$ws_CUL.Cells.Item(129, 13).Value = -997312
$ws_CUL.Cells.Item(129, 14).Value = -14941.4092

$ws_CUL.Cells.Item(131, 8).Value = 4862.8965
$ws_CUL.Cells.Item(131, 10).Value = 1776.9
$ws_CUL.Cells.Item(131, 12).Value = 5330.700000000001
$ws_CUL.Cells.Item(131, 14).Value = -15410.7

$ws_CUL.Cells.Item(137, 8).Value = 12336.477
$ws_CUL.Cells.Item(137, 9).Value = 3410
$ws_CUL.Cells.Item(137, 10).Value = 20451.455
$ws_CUL.Cells.Item(137, 11).Value = 10230
$ws_CUL.Cells.Item(137, 12).Value = 61354.36500000001
$ws_CUL.Cells.Item(137, 13).Value = -5130
$ws_CUL.Cells.Item(137, 14).Value = -71554.36500000001

$ws_CUL.Cells.Item(139, 8).Value = 80167.89999999999
$ws_CUL.Cells.Item(139, 9).Value = 174205.33
$ws_CUL.Cells.Item(139, 10).Value = 3228.182
$ws_CUL.Cells.Item(139, 11).Value = 522615.99
$ws_CUL.Cells.Item(139, 12).Value = 9684.545999999998
$ws_CUL.Cells.Item(139, 13).Value = -517475.99
$ws_CUL.Cells.Item(139, 14).Value = -19964.546

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(122, 8).Value = 1186.3846
$ws_GSM.Cells.Item(122, 9).Value = 1165.8572
$ws_GSM.Cells.Item(122, 10).Value = 1210.3334
$ws_GSM.Cells.Item(122, 11).Value = 3497.5716
$ws_GSM.Cells.Item(122, 12).Value = 3631.0002
$ws_GSM.Cells.Item(122, 13).Value = -1047.5716
$ws_GSM.Cells.Item(122, 14).Value = -8531.0002

$ws_GSM.Cells.Item(132, 8).Value = 55561652
$ws_GSM.Cells.Item(132, 9).Value = 83339790
$ws_GSM.Cells.Item(132, 10).Value = 5371.6665
$ws_GSM.Cells.Item(132, 11).Value = 250019370
$ws_GSM.Cells.Item(132, 12).Value = 16114.9995
$ws_GSM.Cells.Item(132, 13).Value = -250016840
$ws_GSM.Cells.Item(132, 14).Value = -21174.9995

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(50, 8).Value = 10000
$ws_LTW.Cells.Item(50, 9).Value = 10000
$ws_LTW.Cells.Item(50, 11).Value = 10000
$ws_LTW.Cells.Item(50, 13).Value = -9363

$ws_LTW.Cells.Item(94, 8).Value = 33805
$ws_LTW.Cells.Item(94, 10).Value = 33805
$ws_LTW.Cells.Item(94, 12).Value = 33805
$ws_LTW.Cells.Item(94, 14).Value = -35157

$ws_LTW.Cells.Item(132, 8).Value = 3867.2
$ws_LTW.Cells.Item(132, 9).Value = 3161.389
$ws_LTW.Cells.Item(132, 10).Value = 4614.5293
$ws_LTW.Cells.Item(132, 11).Value = 9484.167000000001
$ws_LTW.Cells.Item(132, 12).Value = 13843.5879
$ws_LTW.Cells.Item(132, 13).Value = -6954.167000000001
$ws_LTW.Cells.Item(132, 14).Value = -18903.5879

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(81, 8).Value = 987.625
$ws_WVR.Cells.Item(81, 9).Value = 733.5
$ws_WVR.Cells.Item(81, 10).Value = 1750
$ws_WVR.Cells.Item(81, 11).Value = 1467
$ws_WVR.Cells.Item(81, 12).Value = 3500
$ws_WVR.Cells.Item(81, 13).Value = -406
$ws_WVR.Cells.Item(81, 14).Value = -5622

$ws_WVR.Cells.Item(84, 8).Value = 987.625
$ws_WVR.Cells.Item(84, 9).Value = 733.5
$ws_WVR.Cells.Item(84, 10).Value = 1750
$ws_WVR.Cells.Item(84, 11).Value = 7335
$ws_WVR.Cells.Item(84, 12).Value = 17500
$ws_WVR.Cells.Item(84, 13).Value = -2031
$ws_WVR.Cells.Item(84, 14).Value = -28108

$ws_WVR.Cells.Item(113, 8).Value = 907.6957
$ws_WVR.Cells.Item(113, 9).Value = 630.0833
$ws_WVR.Cells.Item(113, 10).Value = 1210.5454
$ws_WVR.Cells.Item(113, 11).Value = 1890.2499
$ws_WVR.Cells.Item(113, 12).Value = 3631.6362
$ws_WVR.Cells.Item(113, 13).Value = 279.7501
$ws_WVR.Cells.Item(113, 14).Value = -7971.6362

$ws_WVR.Cells.Item(132, 8).Value = 1741364.4
$ws_WVR.Cells.Item(132, 9).Value = 2900088.5
$ws_WVR.Cells.Item(132, 11).Value = 8700265.5
$ws_WVR.Cells.Item(132, 13).Value = -8697735.5
